$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 1.130039812000503
$ws.Range("E3").Value = 1.275815659588881
$ws.Range("F3").Value = 0.5828155900050277
$ws.Range("G3").Value = 0.3008685073931974
$ws.Range("H3").Value = 48.42521619010283

$ws.Range("D9").Value = 1.191894586896629
$ws.Range("E9").Value = 1.514454092434213
$ws.Range("F9").Value = 0.5815323770219396
$ws.Range("G9").Value = 0.2874993982014941
$ws.Range("H9").Value = 51.20941202224161

$ws.Range("D12").Value = 1.772826884560276
$ws.Range("E12").Value = 1.874978778710598
$ws.Range("F12").Value = 0.6130081529339334
$ws.Range("G12").Value = 0.3241854572552367
$ws.Range("H12").Value = 65.42199589408973

$ws.Range("D18").Value = 1.870968211745923
$ws.Range("E18").Value = 1.698016838090766
$ws.Range("F18").Value = 0.6183860409444772
$ws.Range("G18").Value = 0.3289204103846034
$ws.Range("H18").Value = 66.948340593803

$ws.Range("D21").Value = 1.042393725434717
$ws.Range("E21").Value = 0.8462652785780268
$ws.Range("F21").Value = 0.5994023930984035
$ws.Range("G21").Value = 0.2927460172055956
$ws.Range("H21").Value = 42.49750564754882

$ws.Range("D27").Value = 1.359591281829325
$ws.Range("E27").Value = 1.429551824566111
$ws.Range("F27").Value = 0.5965042334550942
$ws.Range("G27").Value = 0.2796140497361679
$ws.Range("H27").Value = 56.12620929339146

$ws.Range("D30").Value = 0.9524985194888282
$ws.Range("E30").Value = 0.6690640918403911
$ws.Range("F30").Value = 0.6078029594341335
$ws.Range("G30").Value = 0.3307531955687637
$ws.Range("H30").Value = 36.18856649138735

$ws.Range("D36").Value = 0.9606868561418768
$ws.Range("E36").Value = 0.6952218388957947
$ws.Range("F36").Value = 0.5980650900149307
$ws.Range("G36").Value = 0.3069677239599218
$ws.Range("H36").Value = 37.74609424586456

$ws.Range("D39").Value = 1.08073914074972
$ws.Range("E39").Value = 0.9193874723821164
$ws.Range("F39").Value = 0.5922852545941724
$ws.Range("G39").Value = 0.3049061961032624
$ws.Range("H39").Value = 45.19627981796808

$ws.Range("D45").Value = 1.65160584235895
$ws.Range("E45").Value = 2.304315262760082
$ws.Range("F45").Value = 0.5997523877657435
$ws.Range("G45").Value = 0.3307587595548684
$ws.Range("H45").Value = 63.68671190281506

$ws.Range("D48").Value = 1.888207575285985
$ws.Range("E48").Value = 1.974674218246942
$ws.Range("F48").Value = 0.6484752350963774
$ws.Range("G48").Value = 0.42350447091781
$ws.Range("H48").Value = 65.65657062369532

$ws.Range("D54").Value = 2.615182453893337
$ws.Range("E54").Value = 3.038460580965362
$ws.Range("F54").Value = 0.6438414075788274
$ws.Range("G54").Value = 0.363440656298077
$ws.Range("H54").Value = 75.38063141176585

$ws.Range("D57").Value = 6.461767059647445
$ws.Range("E57").Value = 11.63974795365123
$ws.Range("F57").Value = 0.6574101129956289
$ws.Range("G57").Value = 0.3861310167055895
$ws.Range("H57").Value = 89.8261558034019

$ws.Range("D63").Value = 8.143928040651026
$ws.Range("E63").Value = 11.29908971185616
$ws.Range("F63").Value = 0.6468857169367002
$ws.Range("G63").Value = 0.3843845969905713
$ws.Range("H63").Value = 92.05683407677816

$ws.Range("D66").Value = 3.429714268763714
$ws.Range("E66").Value = 10.09286962487652
$ws.Range("F66").Value = 0.6318799426006594
$ws.Range("G66").Value = 0.3825753437712935
$ws.Range("H66").Value = 81.57630947990228

$ws.Range("D72").Value = 4.197278605559789
$ws.Range("E72").Value = 9.341076041519006
$ws.Range("F72").Value = 0.6395557151125862
$ws.Range("G72").Value = 0.3621176375219149
$ws.Range("H72").Value = 84.76260989047952

$ws.Range("D75").Value = 2.020363954587721
$ws.Range("E75").Value = 2.293811255464086
$ws.Range("F75").Value = 0.6252891778139488
$ws.Range("G75").Value = 0.3170081591871743
$ws.Range("H75").Value = 69.05066652005547

$ws.Range("D81").Value = 2.351628933613348
$ws.Range("E81").Value = 3.248806248950661
$ws.Range("F81").Value = 0.6225967702997584
$ws.Range("G81").Value = 0.2984945226681516
$ws.Range("H81").Value = 73.52487199827399

$ws.Range("D84").Value = 1.773589740420534
$ws.Range("E84").Value = 2.061381867598989
$ws.Range("F84").Value = 0.5771293380505662
$ws.Range("G84").Value = 0.3016652662707829
$ws.Range("H84").Value = 67.45981751598747

$ws.Range("D90").Value = 2.385951200020108
$ws.Range("E90").Value = 2.925001836010356
$ws.Range("F90").Value = 0.5643341557464877
$ws.Range("G90").Value = 0.2961576268678773
$ws.Range("H90").Value = 76.34762371746196

$ws.Range("D93").Value = 1.537512660735602
$ws.Range("E93").Value = 1.857977131638517
$ws.Range("F93").Value = 0.5975810084176308
$ws.Range("G93").Value = 0.3327908348066964
$ws.Range("H93").Value = 61.13326259494172

$ws.Range("D99").Value = 1.679638470195038
$ws.Range("E99").Value = 1.925028804510991
$ws.Range("F99").Value = 0.5931741671016955
$ws.Range("G99").Value = 0.3149634942528675
$ws.Range("H99").Value = 64.68441407913117

$ws.Range("D102").Value = 1.473442744897767
$ws.Range("E102").Value = 2.024984761337615
$ws.Range("F102").Value = 0.6404060318012496
$ws.Range("G102").Value = 0.3634850751511958
$ws.Range("H102").Value = 56.53675488790824

$ws.Range("D108").Value = 2.095400843123908
$ws.Range("E108").Value = 4.545415804976174
$ws.Range("F108").Value = 0.6410300872521241
$ws.Range("G108").Value = 0.3587792290280774
$ws.Range("H108").Value = 69.40775845558741

$ws.Range("D111").Value = 0.9172207748045611
$ws.Range("E111").Value = 1.310887195018372
$ws.Range("F111").Value = 0.2834395874701757
$ws.Range("G111").Value = 0.06613406133916325
$ws.Range("H111").Value = 69.09799742264121

$ws.Range("D117").Value = 1.004358600562757
$ws.Range("E117").Value = 1.54588599717618
$ws.Range("F117").Value = 0.2811652669671177
$ws.Range("G117").Value = 0.04729203065788602
$ws.Range("H117").Value = 72.00549018950238

$ws.Range("D120").Value = 1.486337451856089
$ws.Range("E120").Value = 1.678731597409714
$ws.Range("F120").Value = 0.4440790461361111
$ws.Range("G120").Value = 0.08247294571698363
$ws.Range("H120").Value = 70.1225959433667

$ws.Range("D126").Value = 2.068229170699706
$ws.Range("E126").Value = 3.319685289326045
$ws.Range("F126").Value = 0.4485041848075516
$ws.Range("G126").Value = 0.09482021137006048
$ws.Range("H126").Value = 78.31457987531346

$ws.Range("D129").Value = 1.81890260086058
$ws.Range("E129").Value = 2.236590111403163
$ws.Range("F129").Value = 0.6133713213769625
$ws.Range("G129").Value = 0.1034113914285109
$ws.Range("H129").Value = 66.2779457741851

$ws.Range("D135").Value = 2.115623043598457
$ws.Range("E135").Value = 3.039212066584865
$ws.Range("F135").Value = 0.6292179154746675
$ws.Range("G135").Value = 0.1117324142731587
$ws.Range("H135").Value = 70.25850529570558

$ws.Range("D138").Value = 3.96523786800304
$ws.Range("E138").Value = 9.074306383649709
$ws.Range("F138").Value = 1.116938443963994
$ws.Range("G138").Value = 0.2320419365983491
$ws.Range("H138").Value = 71.83174172281112

$ws.Range("D144").Value = 4.979707627148833
$ws.Range("E144").Value = 8.730470938522354
$ws.Range("F144").Value = 1.089665349160785
$ws.Range("G144").Value = 0.2093026209643546
$ws.Range("H144").Value = 78.11788500955265

